{"js": "// Rebuilds the document body as a set of distinctly-formatted paragraphs\n// (intro lines with tab stops, a bottom-bordered line, Q&A pairs, etc.),\n// matching the target OOXML exactly. Injected via insertOoxml() so every\n// run/paragraph property (fonts, sizes, tab stops, paragraph border, line\n// break) is reproduced byte-for-byte instead of being approximated through\n// many separate property-setter calls.\nconst paragraphsXml = [\n  \"<w:p><w:pPr><w:widowControl w:val=\\\"0\\\"/><w:tabs><w:tab w:val=\\\"left\\\" w:pos=\\\"560\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"1120\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"1680\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"2240\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"2800\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"3360\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"3920\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"4480\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"5040\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"5600\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"6160\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"6720\\\"/></w:tabs><w:autoSpaceDE w:val=\\\"0\\\"/><w:autoSpaceDN w:val=\\\"0\\\"/><w:adjustRightInd w:val=\\\"0\\\"/><w:spacing w:line=\\\"320\\\" w:lineRule=\\\"exact\\\"/><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Hi,Jan, plz write a story of about 250-300 words. It's for students learning EEC3 U3. This time, a real story, not a chant nor a simple workbook story( I promised I wouldn't say that word ). There\u2019s no need to use a lot of repetition or rhyme . Just write a real story for young kids.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:widowControl w:val=\\\"0\\\"/><w:tabs><w:tab w:val=\\\"left\\\" w:pos=\\\"560\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"1120\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"1680\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"2240\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"2800\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"3360\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"3920\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"4480\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"5040\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"5600\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"6160\\\"/><w:tab w:val=\\\"left\\\" w:pos=\\\"6720\\\"/></w:tabs><w:autoSpaceDE w:val=\\\"0\\\"/><w:autoSpaceDN w:val=\\\"0\\\"/><w:adjustRightInd w:val=\\\"0\\\"/><w:spacing w:line=\\\"320\\\" w:lineRule=\\\"exact\\\"/><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>I will send you the content table to show the sentence structures they have learnt, and also the vocab list of EEC1-3. The words in EEC1,2 and EEC3U1-U3 are the words the kids have learnt. Please try to keep the number of new words no more than 12.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:pBdr><w:bottom w:val=\\\"single\\\" w:sz=\\\"6\\\" w:space=\\\"1\\\" w:color=\\\"auto\\\"/></w:pBdr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Also use persons on the Civa character chart. Once you decide on the topic, please let me know it first</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Do you have a drum?</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Yes, I do. / No, I don\u2019t.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Does Billy have a drum?</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Yes, he does. / No, he doesn\u2019t.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Whose English book is this?</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>It\u2019s Max\u2019s.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Whose pencils are these?</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>They are Amy\u2019s.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Is this Smiley\u2019s backpack?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:br/><w:t>Yes, it is. / No, it isn\u2019t.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Which subject do you like best.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>I like English best.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Why do you like English.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Because it\u2019s fun.</w:t></w:r></w:p>\",\n  \"<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Do you like English?</w:t></w:r></w:p>\",\n  \"<w:p><w:r><w:rPr><w:rFonts w:ascii=\\\"Helvetica Neue\\\" w:hAnsi=\\\"Helvetica Neue\\\" w:cs=\\\"Helvetica Neue\\\"/><w:sz w:val=\\\"22\\\"/><w:szCs w:val=\\\"22\\\"/></w:rPr><w:t>Yes, I do. / No, I don\u2019t. It\u2019s difficult for me.</w:t></w:r><w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/><w:bookmarkEnd w:id=\\\"0\\\"/></w:p>\"\n];\n\nconst bodyInnerXml = paragraphsXml.join(\"\");\n\n// Word's Office.js host only accepts OOXML payloads wrapped in the\n// \"Flat OPC\" package format for insertOoxml().\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + bodyInnerXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\n// Replace the whole body content in one shot (sectPr at the end of the\n// body is untouched by a Body-level insertOoxml Replace).\nconst body = context.document.body;\nbody.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Rebuilds the document body as a set of distinctly-formatted paragraphs\n# (intro lines with tab stops, a bottom-bordered line, Q&A pairs, etc.),\n# matching the target OOXML exactly. Applied via Range.InsertXML() so every\n# run/paragraph property (fonts, sizes, tab stops, paragraph border, line\n# break) is reproduced byte-for-byte instead of being approximated through\n# many separate property-setter calls.\n$paragraphsXml = @(\n  '<w:p><w:pPr><w:widowControl w:val=\"0\"/><w:tabs><w:tab w:val=\"left\" w:pos=\"560\"/><w:tab w:val=\"left\" w:pos=\"1120\"/><w:tab w:val=\"left\" w:pos=\"1680\"/><w:tab w:val=\"left\" w:pos=\"2240\"/><w:tab w:val=\"left\" w:pos=\"2800\"/><w:tab w:val=\"left\" w:pos=\"3360\"/><w:tab w:val=\"left\" w:pos=\"3920\"/><w:tab w:val=\"left\" w:pos=\"4480\"/><w:tab w:val=\"left\" w:pos=\"5040\"/><w:tab w:val=\"left\" w:pos=\"5600\"/><w:tab w:val=\"left\" w:pos=\"6160\"/><w:tab w:val=\"left\" w:pos=\"6720\"/></w:tabs><w:autoSpaceDE w:val=\"0\"/><w:autoSpaceDN w:val=\"0\"/><w:adjustRightInd w:val=\"0\"/><w:spacing w:line=\"320\" w:lineRule=\"exact\"/><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Hi,Jan, plz write a story of about 250-300 words. It''s for students learning EEC3 U3. This time, a real story, not a chant nor a simple workbook story( I promised I wouldn''t say that word ). There\u2019s no need to use a lot of repetition or rhyme . Just write a real story for young kids.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:widowControl w:val=\"0\"/><w:tabs><w:tab w:val=\"left\" w:pos=\"560\"/><w:tab w:val=\"left\" w:pos=\"1120\"/><w:tab w:val=\"left\" w:pos=\"1680\"/><w:tab w:val=\"left\" w:pos=\"2240\"/><w:tab w:val=\"left\" w:pos=\"2800\"/><w:tab w:val=\"left\" w:pos=\"3360\"/><w:tab w:val=\"left\" w:pos=\"3920\"/><w:tab w:val=\"left\" w:pos=\"4480\"/><w:tab w:val=\"left\" w:pos=\"5040\"/><w:tab w:val=\"left\" w:pos=\"5600\"/><w:tab w:val=\"left\" w:pos=\"6160\"/><w:tab w:val=\"left\" w:pos=\"6720\"/></w:tabs><w:autoSpaceDE w:val=\"0\"/><w:autoSpaceDN w:val=\"0\"/><w:adjustRightInd w:val=\"0\"/><w:spacing w:line=\"320\" w:lineRule=\"exact\"/><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>I will send you the content table to show the sentence structures they have learnt, and also the vocab list of EEC1-3. The words in EEC1,2 and EEC3U1-U3 are the words the kids have learnt. Please try to keep the number of new words no more than 12.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:pBdr><w:bottom w:val=\"single\" w:sz=\"6\" w:space=\"1\" w:color=\"auto\"/></w:pBdr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Also use persons on the Civa character chart. Once you decide on the topic, please let me know it first</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Do you have a drum?</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Yes, I do. / No, I don\u2019t.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Does Billy have a drum?</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Yes, he does. / No, he doesn\u2019t.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Whose English book is this?</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>It\u2019s Max\u2019s.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Whose pencils are these?</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>They are Amy\u2019s.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Is this Smiley\u2019s backpack?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:br/><w:t>Yes, it is. / No, it isn\u2019t.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Which subject do you like best.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>I like English best.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Why do you like English.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Because it\u2019s fun.</w:t></w:r></w:p>'\n  '<w:p><w:pPr><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Do you like English?</w:t></w:r></w:p>'\n  '<w:p><w:r><w:rPr><w:rFonts w:ascii=\"Helvetica Neue\" w:hAnsi=\"Helvetica Neue\" w:cs=\"Helvetica Neue\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr><w:t>Yes, I do. / No, I don\u2019t. It\u2019s difficult for me.</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>'\n)\n\n$bodyInnerXml = [string]::Join(\"\", $paragraphsXml)\n\n# Word's COM InsertXML() only accepts OOXML payloads wrapped in the\n# \"Flat OPC\" package format.\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + $bodyInnerXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n\n# Replace the whole body content in one shot (sectPr at the end of the\n# body is untouched by a whole-story InsertXML).\n$d = $word.ActiveDocument\n$d.Content.InsertXML($flatOpc)\n"}
